$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '53.487.77'
$ws.Cells.Item(2, 5).Value = '  -4.90%  '

$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '2.218.32'
$ws.Cells.Item(3, 5).Value = '  -6.36%  '

$ws.Cells.Item(4, 5).Value = '  +0.06%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '483.25'
$ws.Cells.Item(5, 5).Value = '  -3.63%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '125.13'
$ws.Cells.Item(6, 5).Value = '  -3.05%  '

$ws.Cells.Item(7, 5).Value = '  +0.17%  '

$ws.Cells.Item(8, 5).Value = '  -5.15%  '

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '2.223.96'
$ws.Cells.Item(9, 5).Value = '  -6.31%  '

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.0910'
$ws.Cells.Item(10, 5).Value = '  -6.97%  '

$ws.Cells.Item(11, 5).Value = '  -1.29%  '

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '4.66'
$ws.Cells.Item(12, 5).Value = '  -2.21%  '

$ws.Cells.Item(13, 5).Value = '  -2.99%  '

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '2.616.08'
$ws.Cells.Item(14, 5).Value = '  -6.23%  '

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '20.90'
$ws.Cells.Item(15, 5).Value = '  -2.47%  '

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '53.429.95'
$ws.Cells.Item(16, 5).Value = '  -4.91%  '

$ws.Cells.Item(17, 5).Value = '  -3.63%  '

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '2.226.77'
$ws.Cells.Item(18, 5).Value = '  -5.83%  '

$ws.Cells.Item(19, 2).Value = 'Chainlink'
$ws.Cells.Item(19, 3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '9.54'
$ws.Cells.Item(19, 5).Value = '  -4.59%  '

$ws.Cells.Item(20, 2).Value = 'Polkadot'
$ws.Cells.Item(20, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '3.95'
$ws.Cells.Item(20, 5).Value = '  -2.01%  '

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '297.60'
$ws.Cells.Item(21, 5).Value = '  -3.05%  '

$ws.Cells.Item(22, 5).Value = '  -2.67%  '

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '0.999'
$ws.Cells.Item(23, 5).Value = '  -0.17%  '

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '63.25'
$ws.Cells.Item(24, 5).Value = '  -3.12%  '

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '0.998'
$ws.Cells.Item(25, 5).Value = '  -0.44%  '

$ws.Cells.Item(26, 5).Value = '  -2.12%  '

$ws.Cells.Item(27, 5).Value = '  -3.12%  '

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '6.93'
$ws.Cells.Item(28, 5).Value = '  -3.84%  '

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '169.57'
$ws.Cells.Item(29, 5).Value = '  -0.80%  '

$ws.Cells.Item(30, 5).Value = '  -3.95%  '

$ws.Cells.Item(31, 2).Value = 'USDe'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '0.998'
$ws.Cells.Item(31, 5).Value = '  -0.12%  '

$ws.Cells.Item(32, 2).Value = 'PEPE'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '0.0₃0671'
$ws.Cells.Item(32, 5).Value = '  -5.49%  '

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '0.998'
$ws.Cells.Item(33, 5).Value = '  +0.01%  '

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '5.72'
$ws.Cells.Item(34, 5).Value = '  -0.20%  '

$ws.Cells.Item(35, 5).Value = '  -3.87%  '

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '17.38'
$ws.Cells.Item(36, 5).Value = '  -1.05%  '

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '1.14'
$ws.Cells.Item(37, 5).Value = '  -2.38%  '

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.829'
$ws.Cells.Item(38, 5).Value = '  +5.12%  '

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '3.55'
$ws.Cells.Item(39, 5).Value = '  -5.19%  '

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '35.76'
$ws.Cells.Item(40, 5).Value = '  -0.80%  '

$ws.Cells.Item(41, 5).Value = '  -1.73%  '

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '1.35'
$ws.Cells.Item(42, 5).Value = '  -2.02%  '

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '3.26'
$ws.Cells.Item(43, 5).Value = '  -2.55%  '

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '122.38'
$ws.Cells.Item(44, 5).Value = '  -5.60%  '

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '4.60'
$ws.Cells.Item(45, 5).Value = '  -2.39%  '

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.0877'
$ws.Cells.Item(46, 5).Value = '  -2.96%  '

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '0.531'
$ws.Cells.Item(47, 5).Value = '  -5.43%  '

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '228.66'
$ws.Cells.Item(48, 5).Value = '  -4.75%  '

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '0.0467'
$ws.Cells.Item(49, 5).Value = '  -2.65%  '

$ws.Cells.Item(50, 5).Value = '  -3.21%  '

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '15.91'
$ws.Cells.Item(51, 5).Value = '  -5.27%  '
